# Children Tracking: + Multi Child + Add bus info
# Appends 5 new localization rows (bus name, number plate, driver,
# monitor, pick time) to the Sheet1 translation table, reusing the
# existing "lang key" column formatting (Consolas font, vertical-center,
# no wrap) but recolored to a new accent color (FFCE9178), matching the
# style already used for other lang_* key cells (e.g. A101).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Template cell that already carries the "lang key" column style
# (Consolas, vertical-center, no wrap) we want to clone onto column A
# for the newly appended rows.
$styleTemplate = $ws.Range("A101")

$newRows = @(
    @{ Row = 111; Key = "lang_bus_name";      Vi = "Tên tuyến";        En = "Bus's name" },
    @{ Row = 112; Key = "lang_bks";           Vi = "Biển Kiểm Soát";   En = "Number Plate" },
    @{ Row = 113; Key = "lang_driver";        Vi = "Tài Xế";           En = "Driver" },
    @{ Row = 114; Key = "lang_monitor";       Vi = "Giám Sát";         En = "Monitor" },
    @{ Row = 115; Key = "lang_time_bus_pick"; Vi = "Giờ Đón Trả";      En = "Pick Time" }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row

    $cellA = $ws.Cells.Item($rowNum, 1)
    $cellB = $ws.Cells.Item($rowNum, 2)
    $cellC = $ws.Cells.Item($rowNum, 3)

    # Clone the "lang key" cell formatting onto column A, then overwrite
    # the value (Copy brings the style across cleanly, reusing the
    # existing font/style records instead of minting new ones).
    $styleTemplate.Copy($cellA)
    $cellA.Value = $r.Key

    $cellB.Value = $r.Vi
    $cellC.Value = $r.En
}

# New accent font color (FFCE9178) applied only once the rows exist, so
# the single font-table insertion is shared by all 5 key cells.
$keyRange = $ws.Range("A111:A115")
$keyRange.Font.Color = 7901646

$ws.Range("C115").Select()
